# Update Name of Algo
# Updates several numeric prediction values in Sheet1 (RandomForest result data)
# as per the corresponding commit on the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.898600000000007
$ws.Range("E3").Value = 15.8146
$ws.Range("B21").Value = 9.456800000000003
$ws.Range("B23").Value = 9.2232
$ws.Range("E24").Value = 16.5583
$ws.Range("B25").Value = 6.0989
$ws.Range("D27").Value = -8.814300000000005
$ws.Range("D31").Value = -8.782700000000006
$ws.Range("D39").Value = -8.062899999999999
$ws.Range("D48").Value = -7.356499999999998
$ws.Range("D51").Value = -7.591299999999999
$ws.Range("D52").Value = -7.784399999999998
$ws.Range("B53").Value = 5.759599999999998
$ws.Range("D55").Value = -8.318199999999996
$ws.Range("D56").Value = -7.9644
$ws.Range("B57").Value = 4.999399999999997
$ws.Range("D57").Value = -8.541399999999992
$ws.Range("E57").Value = 16.6038
$ws.Range("B59").Value = 4.865999999999999
$ws.Range("E61").Value = 16.4966
$ws.Range("B69").Value = 5.385399999999995
$ws.Range("E70").Value = 17.1085
$ws.Range("D73").Value = -7.720499999999998
$ws.Range("B79").Value = 9.129400000000002
$ws.Range("B83").Value = 5.294999999999996
$ws.Range("E86").Value = 16.71110000000001
$ws.Range("D89").Value = -5.731200000000001
$ws.Range("D90").Value = -8.030799999999999
$ws.Range("B93").Value = 5.731100000000001
$ws.Range("E98").Value = 15.5165
$ws.Range("E100").Value = 16.84460000000001
$ws.Range("E102").Value = 16.68259999999998
